$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-7 from 45243 to 45244
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45244
}
